$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new shop item ("detect item") was added: 水晶球 (Crystal Ball), ItemId 22034013,
# Shelf 3. It is appended as a new row (row 43) of the GameShop data/table "表3"
# right after the last existing row (42).
$ws.Range("A43").Value = 15000042
$ws.Range("B43").Value = 22034013
$ws.Range("C43").Value = 3
$ws.Range("D43").Formula = "=LOOKUP(表3[[#This Row],[ItemId]],[1]其他!`$A:`$A,[1]其他!`$B:`$B)"

# Match the formatting (border/style) used by the rest of the "~Name" column.
$ws.Range("D42").Copy()
$ws.Range("D43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Extend the "表3" table range/autofilter to cover the newly added row.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A3:D43"))

# Leave the view scrolled/selected the way the author last left it.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$null = $ws.Range("C42").Select()
